$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (F column) counts
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F7").Value = 1260
$wsExhibition.Range("F14").Value = 809
$wsExhibition.Range("F19").Value = 1045
$wsExhibition.Range("F23").Value = 397
$wsExhibition.Range("F27").Value = 46

# Sheet "全部类型" (All types) - same events aggregated, update matching rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 1260
$wsAll.Range("F21").Value = 809
$wsAll.Range("F26").Value = 1045
$wsAll.Range("F32").Value = 397
$wsAll.Range("F42").Value = 46
